# Auto-generated script to apply profit-sheet corrections
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 452.07144
$ws.Range("I28").Value = 452.07144
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 452.07144
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 32.92856
$ws.Range("H51").Value = 64181988
$ws.Range("I51").Value = 111447784
$ws.Range("J51").Value = 50002250
$ws.Range("K51").Value = 111447784
$ws.Range("L51").Value = 50002250
$ws.Range("M51").Value = -111447300
$ws.Range("N51").Value = -50003218
$ws.Range("H64").Value = 6408.6875
$ws.Range("I64").Value = 7141.2856
$ws.Range("J64").Value = 5838.8887
$ws.Range("K64").Value = 7141.2856
$ws.Range("L64").Value = 5838.8887
$ws.Range("M64").Value = -6893.2856
$ws.Range("N64").Value = -6334.8887
$ws.Range("H67").Value = 6408.6875
$ws.Range("I67").Value = 7141.2856
$ws.Range("J67").Value = 5838.8887
$ws.Range("K67").Value = 7141.2856
$ws.Range("L67").Value = 5838.8887
$ws.Range("M67").Value = -6283.2856
$ws.Range("N67").Value = -7554.8887
$ws.Range("H74").Value = 7526.6
$ws.Range("I74").Value = 9206.817999999999
$ws.Range("K74").Value = 9206.817999999999
$ws.Range("M74").Value = -8270.817999999999
$ws.Range("H76").Value = 3465.6667
$ws.Range("I76").Value = 3465.6667
$ws.Range("K76").Value = 3465.6667
$ws.Range("M76").Value = -3150.6667
$ws.Range("H77").Value = 7526.6
$ws.Range("I77").Value = 9206.817999999999
$ws.Range("K77").Value = 46034.09
$ws.Range("M77").Value = -41354.09
$ws.Range("H79").Value = 3465.6667
$ws.Range("I79").Value = 3465.6667
$ws.Range("K79").Value = 3465.6667
$ws.Range("M79").Value = -2373.6667
$ws.Range("H95").Value = 16000
$ws.Range("J95").Value = 16000
$ws.Range("L95").Value = 16000
$ws.Range("N95").Value = -21492
$ws.Range("H98").Value = 3731.2
$ws.Range("I98").Value = 3671.75
$ws.Range("K98").Value = 3671.75
$ws.Range("M98").Value = -2173.75
$ws.Range("H103").Value = 585.75
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("H105").Value = 37174.25
$ws.Range("J105").Value = 37174.25
$ws.Range("L105").Value = 37174.25
$ws.Range("N105").Value = -44162.25
$ws.Range("H111").Value = 2636.2354
$ws.Range("I111").Value = 2678
$ws.Range("J111").Value = 2536
$ws.Range("K111").Value = 8034
$ws.Range("L111").Value = 7608
$ws.Range("M111").Value = -4967
$ws.Range("N111").Value = -13742
$ws.Range("H122").Value = 3731.2
$ws.Range("I122").Value = 3671.75
$ws.Range("K122").Value = 11015.25
$ws.Range("M122").Value = -8565.25
$ws.Range("H132").Value = 20038.59
$ws.Range("I132").Value = 25924.691
$ws.Range("K132").Value = 77774.073
$ws.Range("M132").Value = -75244.073
$ws.Range("N28").ClearContents()
$ws.Range("M103").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1381.7273
$ws.Range("I2").Value = 1381.7273
$ws.Range("K2").Value = 1381.7273
$ws.Range("M2").Value = -1268.7273
$ws.Range("H14").Value = 1748
$ws.Range("I14").Value = 1748
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1748
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1573
$ws.Range("H45").Value = 75460.71000000001
$ws.Range("I45").Value = 103395.5
$ws.Range("K45").Value = 103395.5
$ws.Range("M45").Value = -103018.5
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("H104").Value = 68500
$ws.Range("J104").Value = 68500
$ws.Range("L104").Value = 68500
$ws.Range("N104").Value = -75488
$ws.Range("H107").Value = 44113
$ws.Range("J107").Value = 44113
$ws.Range("L107").Value = 44113
$ws.Range("N107").Value = -51793
$ws.Range("H116").Value = 1381.7273
$ws.Range("I116").Value = 1381.7273
$ws.Range("K116").Value = 1381.7273
$ws.Range("M116").Value = 912.2727
$ws.Range("H132").Value = 2321.923
$ws.Range("I132").Value = 1467.5264
$ws.Range("J132").Value = 4641
$ws.Range("K132").Value = 4402.5792
$ws.Range("L132").Value = 13923
$ws.Range("M132").Value = -1872.5792
$ws.Range("N132").Value = -18983
$ws.Range("N14").ClearContents()
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1381.7273
$ws.Range("I3").Value = 1381.7273
$ws.Range("K3").Value = 1381.7273
$ws.Range("M3").Value = -1267.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2673.6143
$ws.Range("J31").Value = 3059.9424
$ws.Range("L31").Value = 3059.9424
$ws.Range("N31").Value = -3649.9424
$ws.Range("H34").Value = 2673.6143
$ws.Range("J34").Value = 3059.9424
$ws.Range("L34").Value = 3059.9424
$ws.Range("N34").Value = -3463.9424
$ws.Range("H92").Value = 58249.5
$ws.Range("J92").Value = 58249.5
$ws.Range("L92").Value = 58249.5
$ws.Range("N92").Value = -63241.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2614.9167
$ws.Range("I48").Value = 1444.4445
$ws.Range("K48").Value = 4333.333500000001
$ws.Range("M48").Value = -4083.333500000001
$ws.Range("H97").Value = 461
$ws.Range("I97").Value = 464.66666
$ws.Range("J97").Value = 450
$ws.Range("K97").Value = 1393.99998
$ws.Range("L97").Value = 1350
$ws.Range("M97").Value = -897.9999800000001
$ws.Range("N97").Value = -2342
$ws.Range("H98").Value = 494.7143
$ws.Range("I98").Value = 355
$ws.Range("J98").Value = 844
$ws.Range("K98").Value = 1065
$ws.Range("L98").Value = 2532
$ws.Range("M98").Value = 433
$ws.Range("N98").Value = -5528

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 15277
$ws.Range("J33").Value = 15277
$ws.Range("L33").Value = 15277
$ws.Range("N33").Value = -15781
$ws.Range("H36").Value = 12999
$ws.Range("J36").Value = 12999
$ws.Range("L36").Value = 12999
$ws.Range("N36").Value = -13969
$ws.Range("H70").Value = 11632.863
$ws.Range("I70").Value = 10022.315
$ws.Range("K70").Value = 10022.315
$ws.Range("M70").Value = -9752.315000000001
$ws.Range("H73").Value = 11632.863
$ws.Range("I73").Value = 10022.315
$ws.Range("K73").Value = 10022.315
$ws.Range("M73").Value = -9086.315000000001
$ws.Range("H80").Value = 134414.06
$ws.Range("I80").Value = 150735.78
$ws.Range("J80").Value = 67495
$ws.Range("K80").Value = 150735.78
$ws.Range("L80").Value = 67495
$ws.Range("M80").Value = -149737.78
$ws.Range("N80").Value = -69491
$ws.Range("H83").Value = 134414.06
$ws.Range("I83").Value = 150735.78
$ws.Range("J83").Value = 67495
$ws.Range("K83").Value = 753678.9
$ws.Range("L83").Value = 337475
$ws.Range("M83").Value = -748686.9
$ws.Range("N83").Value = -347459
$ws.Range("H132").Value = 25061764
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11349.647
$ws.Range("J7").Value = 18932.625
$ws.Range("L7").Value = 18932.625
$ws.Range("N7").Value = -19156.625
$ws.Range("H23").Value = 515817.7
$ws.Range("I23").Value = 515817.7
$ws.Range("K23").Value = 515817.7
$ws.Range("M23").Value = -515587.7
$ws.Range("H40").Value = 3090813.5
$ws.Range("I40").Value = 4633579.5
$ws.Range("K40").Value = 4633579.5
$ws.Range("M40").Value = -4633443.5
$ws.Range("H61").Value = 2991.0938
$ws.Range("I61").Value = 2860.2964
$ws.Range("J61").Value = 3697.4
$ws.Range("K61").Value = 2860.2964
$ws.Range("L61").Value = 3697.4
$ws.Range("M61").Value = -2658.2964
$ws.Range("N61").Value = -4101.4
$ws.Range("H93").Value = 3384.6155
$ws.Range("I93").Value = 2275
$ws.Range("J93").Value = 5160
$ws.Range("K93").Value = 2275
$ws.Range("L93").Value = 5160
$ws.Range("M93").Value = -1027
$ws.Range("N93").Value = -7656
$ws.Range("H94").Value = 55999.8
$ws.Range("H100").Value = 2084.1538
$ws.Range("I100").Value = 1843.8889
$ws.Range("K100").Value = 1843.8889
$ws.Range("M100").Value = -1302.8889
$ws.Range("H113").Value = 2991.0938
$ws.Range("I113").Value = 2860.2964
$ws.Range("J113").Value = 3697.4
$ws.Range("K113").Value = 2860.2964
$ws.Range("L113").Value = 3697.4
$ws.Range("M113").Value = -690.2964000000002
$ws.Range("N113").Value = -8037.4
$ws.Range("H126").Value = 11349.647
$ws.Range("J126").Value = 18932.625
$ws.Range("L126").Value = 56797.875
$ws.Range("N126").Value = -61737.875
$ws.Range("H132").Value = 6631.4414
$ws.Range("I132").Value = 2894.7
$ws.Range("J132").Value = 11969.643
$ws.Range("K132").Value = 8684.099999999999
$ws.Range("L132").Value = 35908.929
$ws.Range("M132").Value = -6154.099999999999
$ws.Range("N132").Value = -40968.929
$ws.Range("H136").Value = 5229.231
$ws.Range("I136").Value = 4856.4287
$ws.Range("K136").Value = 14569.2861
$ws.Range("M136").Value = -12019.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 14572
$ws.Range("J97").Value = 14572
$ws.Range("L97").Value = 14572
$ws.Range("N97").Value = -16554
$ws.Range("H132").Value = 1689.8206
$ws.Range("I132").Value = 1674.0454
$ws.Range("J132").Value = 1710.2354
$ws.Range("K132").Value = 5022.1362
$ws.Range("L132").Value = 5130.706200000001
$ws.Range("M132").Value = -2492.1362
$ws.Range("N132").Value = -10190.7062
